$d = $word.ActiveDocument

# Locate the paragraph with the exact text "Berikut adalah langkah-langkah Maintenance :"
# (there are similarly worded paragraphs for Alpha/Beta testing in this document, so we
# must match the full text exactly to avoid touching the wrong one).

$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text.Trim()
    if ($t -eq "Berikut adalah langkah-langkah Maintenance :") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $r = $target.Range
    $r.Collapse(1)

    # New paragraph to insert right before the target paragraph.
    $newParaXml = '<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr>' + `
        '<w:r><w:t>Maintenance website adalah kegiatan pemeliharaan website yang bertujuan untuk merawat website agar tetap berada pada performa yang baik, ter</w:t></w:r>' + `
        '<w:r><w:t>-</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">update, dan terhindar dari berbagai permasalahan yang dapat </w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">merusak atau </w:t></w:r>' + `
        '<w:r><w:t>merugikan</w:t></w:r>' + `
        '<w:r><w:t>.</w:t></w:r>' + `
        '</w:p>'

    # Re-emit the original target paragraph, now with a first-line indent added to its pPr.
    $targetParaXml = '<w:p w14:paraId="1170E135" w14:textId="522BDD76" w:rsidR="00C97093" w:rsidRDefault="00F17E0E" w:rsidP="00C97093">' + `
        '<w:pPr><w:ind w:firstLine="720"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>' + `
        '<w:r w:rsidRPr="006B34E5"><w:t xml:space="preserve">Berikut adalah langkah-langkah </w:t></w:r>' + `
        '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Maintenance :</w:t></w:r>' + `
        '</w:p>'

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $newParaXml + $targetParaXml + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$r.InsertXML($xml)
}
